$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 2.242386666666667
$ws.Cells.Item(2,8).Value = 6.72716
$ws.Cells.Item(2,9).Value = 0.04442500453715972
$ws.Cells.Item(2,10).Value = 0.04442500453715972
$ws.Cells.Item(2,13).Value = 247.0944516666667
$ws.Cells.Item(2,14).Value = 741.283355
$ws.Cells.Item(2,15).Value = 0.8050739182622993
$ws.Cells.Item(2,16).Value = 0.8050739182622993
$ws.Cells.Item(2,17).Value = 554.0813038246445
$ws.Cells.Item(2,18).Value = 4986.7317344218
$ws.Cells.Item(2,19).Value = 0.03576541247155161
$ws.Cells.Item(2,20).Value = 0.0357654124715516
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 2.242386666666667
$ws.Cells.Item(3,8).Value = 6.72716
$ws.Cells.Item(3,9).Value = 0.04442500453715972
$ws.Cells.Item(3,10).Value = 0.04442500453715972
$ws.Cells.Item(3,15).Value = 0.1379009747488701
$ws.Cells.Item(3,16).Value = 0.13790097474887
$ws.Cells.Item(3,17).Value = 94.90849244311111
$ws.Cells.Item(3,18).Value = 854.176431988
$ws.Cells.Item(3,19).Value = 0.006126251428897301
$ws.Cells.Item(3,20).Value = 0.006126251428897299
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 2.242386666666667
$ws.Cells.Item(4,8).Value = 6.72716
$ws.Cells.Item(4,9).Value = 0.04442500453715972
$ws.Cells.Item(4,10).Value = 0.04442500453715972
$ws.Cells.Item(4,13).Value = 11.590146
$ws.Cells.Item(4,14).Value = 34.770438
$ws.Cells.Item(4,15).Value = 0.03776258103132013
$ws.Cells.Item(4,16).Value = 0.03776258103132013
$ws.Cells.Item(4,17).Value = 25.98958885512
$ws.Cells.Item(4,18).Value = 233.90629969608
$ws.Cells.Item(4,19).Value = 0.001677602833651259
$ws.Cells.Item(4,20).Value = 0.001677602833651258
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 2.242386666666667
$ws.Cells.Item(5,8).Value = 6.72716
$ws.Cells.Item(5,9).Value = 0.04442500453715972
$ws.Cells.Item(5,10).Value = 0.04442500453715972
$ws.Cells.Item(5,13).Value = 5.912082333333333
$ws.Cells.Item(5,14).Value = 17.736247
$ws.Cells.Item(5,15).Value = 0.01926252595751047
$ws.Cells.Item(5,16).Value = 0.01926252595751047
$ws.Cells.Item(5,17).Value = 13.25717459650222
$ws.Cells.Item(5,18).Value = 119.31457136852
$ws.Cells.Item(5,19).Value = 0.0008557378030595597
$ws.Cells.Item(5,20).Value = 0.0008557378030595596
$ws.Cells.Item(6,9).Value = 0.4052409520727612
$ws.Cells.Item(6,10).Value = 0.4052409520727612
$ws.Cells.Item(6,13).Value = 247.0944516666667
$ws.Cells.Item(6,14).Value = 741.283355
$ws.Cells.Item(6,15).Value = 0.8050739182622993
$ws.Cells.Item(6,16).Value = 0.8050739182622993
$ws.Cells.Item(6,17).Value = 5054.280521227637
$ws.Cells.Item(6,18).Value = 45488.52469104873
$ws.Cells.Item(6,19).Value = 0.3262489211255625
$ws.Cells.Item(6,20).Value = 0.3262489211255625
$ws.Cells.Item(7,9).Value = 0.4052409520727612
$ws.Cells.Item(7,10).Value = 0.4052409520727612
$ws.Cells.Item(7,15).Value = 0.1379009747488701
$ws.Cells.Item(7,16).Value = 0.13790097474887
$ws.Cells.Item(7,19).Value = 0.05588312229899391
$ws.Cells.Item(7,20).Value = 0.05588312229899389
$ws.Cells.Item(8,9).Value = 0.4052409520727612
$ws.Cells.Item(8,10).Value = 0.4052409520727612
$ws.Cells.Item(8,13).Value = 11.590146
$ws.Cells.Item(8,14).Value = 34.770438
$ws.Cells.Item(8,15).Value = 0.03776258103132013
$ws.Cells.Item(8,16).Value = 0.03776258103132013
$ws.Cells.Item(8,17).Value = 237.074724951774
$ws.Cells.Item(8,18).Value = 2133.672524565966
$ws.Cells.Item(8,19).Value = 0.01530294428985696
$ws.Cells.Item(8,20).Value = 0.01530294428985696
$ws.Cells.Item(9,9).Value = 0.4052409520727612
$ws.Cells.Item(9,10).Value = 0.4052409520727612
$ws.Cells.Item(9,13).Value = 5.912082333333333
$ws.Cells.Item(9,14).Value = 17.736247
$ws.Cells.Item(9,15).Value = 0.01926252595751047
$ws.Cells.Item(9,16).Value = 0.01926252595751047
$ws.Cells.Item(9,17).Value = 120.9307711108421
$ws.Cells.Item(9,18).Value = 1088.376939997579
$ws.Cells.Item(9,19).Value = 0.00780596435834782
$ws.Cells.Item(9,20).Value = 0.00780596435834782
$ws.Cells.Item(10,7).Value = 27.778539
$ws.Cells.Item(10,8).Value = 83.335617
$ws.Cells.Item(10,9).Value = 0.5503340433900792
$ws.Cells.Item(10,10).Value = 0.5503340433900791
$ws.Cells.Item(10,13).Value = 247.0944516666667
$ws.Cells.Item(10,14).Value = 741.283355
$ws.Cells.Item(10,15).Value = 0.8050739182622993
$ws.Cells.Item(10,16).Value = 0.8050739182622993
$ws.Cells.Item(10,17).Value = 6863.922862306115
$ws.Cells.Item(10,18).Value = 61775.30576075504
$ws.Cells.Item(10,19).Value = 0.4430595846651853
$ws.Cells.Item(10,20).Value = 0.4430595846651852
$ws.Cells.Item(11,7).Value = 27.778539
$ws.Cells.Item(11,8).Value = 83.335617
$ws.Cells.Item(11,9).Value = 0.5503340433900792
$ws.Cells.Item(11,10).Value = 0.5503340433900791
$ws.Cells.Item(11,15).Value = 0.1379009747488701
$ws.Cells.Item(11,16).Value = 0.13790097474887
$ws.Cells.Item(11,17).Value = 1175.7201815159
$ws.Cells.Item(11,18).Value = 10581.4816336431
$ws.Cells.Item(11,19).Value = 0.07589160102097887
$ws.Cells.Item(11,20).Value = 0.07589160102097883
$ws.Cells.Item(12,7).Value = 27.778539
$ws.Cells.Item(12,8).Value = 83.335617
$ws.Cells.Item(12,9).Value = 0.5503340433900792
$ws.Cells.Item(12,10).Value = 0.5503340433900791
$ws.Cells.Item(12,13).Value = 11.590146
$ws.Cells.Item(12,14).Value = 34.770438
$ws.Cells.Item(12,15).Value = 0.03776258103132013
$ws.Cells.Item(12,16).Value = 0.03776258103132013
$ws.Cells.Item(12,17).Value = 321.9573226766939
$ws.Cells.Item(12,18).Value = 2897.615904090246
$ws.Cells.Item(12,19).Value = 0.02078203390781191
$ws.Cells.Item(12,20).Value = 0.02078203390781191
$ws.Cells.Item(13,7).Value = 27.778539
$ws.Cells.Item(13,8).Value = 83.335617
$ws.Cells.Item(13,9).Value = 0.5503340433900792
$ws.Cells.Item(13,10).Value = 0.5503340433900791
$ws.Cells.Item(13,13).Value = 5.912082333333333
$ws.Cells.Item(13,14).Value = 17.736247
$ws.Cells.Item(13,15).Value = 0.01926252595751047
$ws.Cells.Item(13,16).Value = 0.01926252595751047
$ws.Cells.Item(13,17).Value = 164.229009667711
$ws.Cells.Item(13,18).Value = 1478.061087009399
$ws.Cells.Item(13,19).Value = 0.0106008237961031
$ws.Cells.Item(13,20).Value = 0.01060082379610309
